{"js": "// Update the phone number in the resume's contact line:\n//   +1-(206)-327-8537  ->  +1-(541)-604-2147\n// The contact line (\"Duc-Huy, DO | +1-(206)-327-8537 | DUCHUYdo.DDH@GMAIL.COM\")\n// lives inside a Word content control (a data-bound \"Your Name\" sdt), so we\n// locate it through the paragraph's parent content control and replace its\n// text in place (search-range mutation does not reliably land inside sdt\n// content in this host, but ContentControl.insertText does).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst oldPhone = \"(206)-327-8537\";\nconst newPhone = \"(541)-604-2147\";\n\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(oldPhone) !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!targetParagraph) {\n  throw new Error(\"Could not find the paragraph containing the phone number.\");\n}\n\nconst cc = targetParagraph.parentContentControlOrNullObject;\ncc.load(\"text,isNullObject\");\nawait context.sync();\n\nif (!cc.isNullObject && cc.text.indexOf(oldPhone) !== -1) {\n  // The contact line is wrapped in a content control: replace the whole\n  // control's text in place, preserving the surrounding run formatting.\n  const updatedText = cc.text.split(oldPhone).join(newPhone);\n  cc.insertText(updatedText, Word.InsertLocation.replace);\n} else {\n  // Fallback: no content control wrapping the text \u2014 replace directly via\n  // the paragraph itself.\n  const updatedText = targetParagraph.text.split(oldPhone).join(newPhone);\n  targetParagraph.insertText(updatedText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the phone number in the resume's contact line:\n#   +1-(206)-327-8537  ->  +1-(541)-604-2147\n# The contact line (\"Duc-Huy, DO | +1-(206)-327-8537 | DUCHUYdo.DDH@GMAIL.COM\")\n# lives inside a Word content control (a data-bound \"Your Name\" sdt), so we\n# update it through the ContentControls collection rather than a plain\n# Range/Find replace, which does not reliably persist when the matched text\n# sits inside content-control-wrapped content.\n\n$d = $word.ActiveDocument\n\n$oldPhone = \"(206)-327-8537\"\n$newPhone = \"(541)-604-2147\"\n\n$updated = $false\n\nfor ($i = 1; $i -le $d.ContentControls.Count; $i++) {\n    $cc = $d.ContentControls.Item($i)\n    if ($cc.Range.Text -like \"*$oldPhone*\") {\n        $cc.Range.Text = $cc.Range.Text.Replace($oldPhone, $newPhone)\n        $updated = $true\n    }\n}\n\nif (-not $updated) {\n    # Fallback: no content control matched \u2014 fall back to a direct Find/Range replace.\n    $rng = $d.Content\n    $f = $rng.Find\n    $f.Text = $oldPhone\n    $f.Forward = $true\n    $f.Wrap = 1\n    if ($f.Execute()) {\n        $rng.Text = $newPhone\n    }\n}\n"}
